# Daily attendance processing - 2026-01-27 06:09:16
# Swap the order of names in column G ("Recorded By") from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact value is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
